$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '54.623.82'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.71%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.279.92'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '505.93'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.65%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.86'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.30%  '

$ws.Range("E7").Value = '  -0.23%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.527'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.15%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.295.60'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.45%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0970'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.87%  '

$ws.Range("E11").Value = '  +1.51%  '

$ws.Range("E12").Value = '  +2.37%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.91'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.07%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.47'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.94%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.686.47'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.12%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '54.669.75'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.88%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.297.19'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.18%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.39'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.69%  '

$ws.Range("E20").Value = '  +0.97%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '307.35'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.20%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.50'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.60%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.09%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.35'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.33%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.996'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.39%  '

$ws.Range("E26").Value = '  -0.51%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.47'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.18%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '171.55'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.96%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.09'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.21%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0704'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.71%  '

$ws.Range("E31").Value = '  +0.92%  '

$ws.Range("E32").Value = '  +3.96%  '

$ws.Range("E33").Value = '  -0.02%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.96'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.24%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.994'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.27%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.910'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.52%  '

$ws.Range("E37").Value = '  +0.57%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.80'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.72%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.58'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.40%  '

$ws.Range("E40").Value = '  +0.30%  '

$ws.Range("E41").Value = '  +1.08%  '

$ws.Range("E42").Value = '  +0.65%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '129.45'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.59%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.97'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.47%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '250.76'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.49%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0498'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.14%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0908'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.43%  '

$ws.Range("E48").Value = '  +0.84%  '

$ws.Range("E49").Value = '  +0.40%  '

$ws.Range("E50").Value = '  +0.70%  '

$ws.Range("E51").Value = '  +0.36%  '
